# 06_Relazionale_v2.xlsx - "ER e ER ristr, Relazionale"
#
# Updates the NOT NULL / UNIQUE-NOT NULL legend markers on the
# "Relazionale" worksheet (Foglio1):
#   - H7 / K7  : "UNIQUE - NOT NULL" -> "NOT NULL"
#   - F15 / G15: "NOT NULL" -> (cleared)
#   - M31 / N31: (empty)    -> "NOT NULL"
# Also refreshes the view (zoom/scroll/selection) to match where the
# author ended up working last.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# --- content edits -------------------------------------------------
$ws.Range("H7").Value = "NOT NULL"
$ws.Range("K7").Value = "NOT NULL"

$ws.Range("F15").Value = ""
$ws.Range("G15").Value = ""

$ws.Range("M31").Value = "NOT NULL"
$ws.Range("N31").Value = "NOT NULL"

# --- view / window state --------------------------------------------
$excel.ActiveWindow.Zoom = 85

[void]$ws.Range("D7").Select()
[void]$ws.Range("M31").Select()
